$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) for rows 2-4
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum(B:E)

$data = @{
    2 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    3 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    4 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 3.537761648806719; E = 0.4942365360607697; G = 8.974608811992548 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
